# Webdata_TestData.xlsx edit: "Added the Product for the Parent Company"
#
# 1. Rename the sanity/parent-company related text that was used for the old
#    "Automated Sanity2" run to the new "Web Data Sanity" naming, everywhere
#    it appears in the workbook (shared string, so touches every sheet that
#    references the parent company name / its derived child+reseller names).
# 2. Move the active selection on a couple of sheets (LoginData, ConfigAccType)
#    and move the active tab/selection over to the AddProduct sheet, where the
#    new product for the parent company was added.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Automated Sanity2", "Web Data Sanity")
}

$wsLogin = $wb.Worksheets.Item("LoginData")
$wsLogin.Activate()
[void]$wsLogin.Range("E8").Select()

$wsConfigAccType = $wb.Worksheets.Item("ConfigAccType")
$wsConfigAccType.Activate()
[void]$wsConfigAccType.Range("A27").Select()

$wsAddProduct = $wb.Worksheets.Item("AddProduct")
$wsAddProduct.Activate()
[void]$wsAddProduct.Range("S1").Select()

Write-Output "Renamed parent company to 'Web Data Sanity' and set AddProduct as the active sheet."
